$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New vintage column BH: "Agosto.2021" ------------------------------
$ws.Cells.Item(1, 60).Value = "Agosto.2021"

# Copy the header formatting (bold, centered, bordered) from the previous
# vintage header cell instead of assigning the .Style COM object directly.
$ws.Cells.Item(1, 59).Copy() | Out-Null
$ws.Cells.Item(1, 60).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Carry the last reported figure for each series forward into the new
# column (matches how every previous vintage column repeats the latest
# known value until a newer data point replaces it).
for ($r = 2; $r -le 73; $r++) {
    $prev = $ws.Cells.Item($r, 59).Value2
    if ($prev -ne $null) {
        $ws.Cells.Item($r, 60).Value2 = $prev
    }
}

# Row 74 (01-01-2021) receives a revised figure in this vintage.
$ws.Cells.Item(74, 60).Value2 = 10937

# --- New row 75: the next quarter, 01-04-2021 ---------------------------
# Use a text formula + paste-values round trip so the date-like label is
# stored as literal text (same as every other "Serie" label in column A)
# instead of being auto-converted into a date serial number.
$ws.Cells.Item(75, 1).Formula = '="01-04-2021"'
$ws.Cells.Item(75, 1).Copy() | Out-Null
$ws.Cells.Item(75, 1).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(75, 60).Value2 = 10700
